$d = $word.ActiveDocument

# The document currently carries a "_GoBack" bookmark (marking the last edit location)
# around an empty paragraph near the end of the document. That bookmark is about to be
# re-anchored inside the "Descriptions" field note below, so remove the old one first -
# bookmark names must stay unique.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

function Insert-PlainFieldNote {
    param(
        [int]$ParaIndex,
        [string]$LabelText,
        [bool]$LabelPreserveSpace,
        [string]$AppendXml
    )

    $p = $d.Paragraphs.Item($ParaIndex)
    $pPrXml = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>'
    if ($LabelPreserveSpace) {
        $spaceAttr = ' xml:space="preserve"'
    } else {
        $spaceAttr = ''
    }
    $labelRun = '<w:r><w:t' + $spaceAttr + '>' + $LabelText + '</w:t></w:r>'

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
        $pPrXml + $labelRun + $AppendXml + `
        '</w:p></w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $p.Range.InsertXML($xml) | Out-Null
}

# "Title" paragraph (#8): add "(Letter and number and symbos(sanitize))" with spell-check markers
$appendXml = '<w:r><w:t xml:space="preserve"> (Letter and number and </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>symbos</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>(sanitize))</w:t></w:r>'
Insert-PlainFieldNote 8 "Title" $false $appendXml

# "Due date" paragraph (#9): add " (date format)"
$appendXml = '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>(date format)</w:t></w:r>'
Insert-PlainFieldNote 9 "Due date" $false $appendXml

# "Descriptions " paragraph (#10): add "(Letter and number and symbos(sanitize))" with the _GoBack bookmark
# wrapped around the spell-checked word, matching the diff.
$appendXml = '<w:r><w:t xml:space="preserve">(Letter and number and </w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>symbos</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>(sanitize))</w:t></w:r>' + `
    '<w:bookmarkEnd w:id="0"/>'
Insert-PlainFieldNote 10 "Descriptions " $true $appendXml

# "Status " paragraph (#11): add "(static)"
$appendXml = '<w:r><w:t>(static)</w:t></w:r>'
Insert-PlainFieldNote 11 "Status " $true $appendXml

# "Date creation " paragraph (#12): add "(date format)"
$appendXml = '<w:r><w:t>(date format)</w:t></w:r>'
Insert-PlainFieldNote 12 "Date creation " $true $appendXml

# "Priority Number " paragraph (#13): add "(static) "
$appendXml = '<w:r><w:t xml:space="preserve">(static) </w:t></w:r>'
Insert-PlainFieldNote 13 "Priority Number " $true $appendXml
